$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 278 ("Fruta / hortaliza, semanal").
# This shifts the previously-existing rows 278-415 down to 279-416.
$ws.Rows.Item(278).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(278, 1).Value  = 4
$ws.Cells.Item(278, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(278, 3).Value  = "Los Lagos"
$ws.Cells.Item(278, 4).Value  = 44992
$ws.Cells.Item(278, 5).Value  = 10
$ws.Cells.Item(278, 6).Value  = 100112037
$ws.Cells.Item(278, 7).Value  = "Cebollín"
$ws.Cells.Item(278, 8).Value  = "Sin especificar"
$ws.Cells.Item(278, 9).Value  = "Primera"
$ws.Cells.Item(278, 10).Value = 180
$ws.Cells.Item(278, 11).Value = 6500
$ws.Cells.Item(278, 12).Value = 7500
$ws.Cells.Item(278, 13).Value = 7000
$ws.Cells.Item(278, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(278, 15).Value = "Región Metropolitana"
$ws.Cells.Item(278, 16).Value = 194
$ws.Cells.Item(278, 17).Value = 36
$ws.Cells.Item(278, 18).Value = "Hortaliza"
